# Auto-generated script applying market-price / profit value updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 187.54546
$ws.Range("I11").Value = 187.54546
$ws.Range("K11").Value = 187.54546
$ws.Range("M11").Value = -47.54545999999999
$ws.Range("H17").Value = 4214.846
$ws.Range("J17").Value = 4214.846
$ws.Range("L17").Value = 12644.538
$ws.Range("N17").Value = -12980.538
$ws.Range("H51").Value = 12057.546
$ws.Range("I51").Value = 17200.2
$ws.Range("J51").Value = 10545
$ws.Range("K51").Value = 17200.2
$ws.Range("L51").Value = 10545
$ws.Range("M51").Value = -16716.2
$ws.Range("N51").Value = -11513
$ws.Range("H55").Value = 391.5
$ws.Range("J55").Value = 799
$ws.Range("L55").Value = 799
$ws.Range("N55").Value = -1227
$ws.Range("H62").Value = 5120.6665
$ws.Range("I62").Value = 3793.1538
$ws.Range("J62").Value = 13749.5
$ws.Range("K62").Value = 3793.1538
$ws.Range("L62").Value = 13749.5
$ws.Range("M62").Value = -3169.1538
$ws.Range("N62").Value = -14997.5
$ws.Range("H65").Value = 5120.6665
$ws.Range("I65").Value = 3793.1538
$ws.Range("J65").Value = 13749.5
$ws.Range("K65").Value = 18965.769
$ws.Range("L65").Value = 68747.5
$ws.Range("M65").Value = -15845.769
$ws.Range("N65").Value = -74987.5
$ws.Range("H74").Value = 14321.643
$ws.Range("I74").Value = 16900.363
$ws.Range("J74").Value = 4866.3335
$ws.Range("K74").Value = 16900.363
$ws.Range("L74").Value = 4866.3335
$ws.Range("M74").Value = -15964.363
$ws.Range("N74").Value = -6738.3335
$ws.Range("H77").Value = 14321.643
$ws.Range("I77").Value = 16900.363
$ws.Range("J77").Value = 4866.3335
$ws.Range("K77").Value = 84501.815
$ws.Range("L77").Value = 24331.6675
$ws.Range("M77").Value = -79821.815
$ws.Range("N77").Value = -33691.6675
$ws.Range("H98").Value = 889.4211
$ws.Range("I98").Value = 844.3889
$ws.Range("K98").Value = 844.3889
$ws.Range("M98").Value = 653.6111
$ws.Range("H106").Value = 88004910
$ws.Range("I106").Value = 220005000
$ws.Range("J106").Value = 4852.6665
$ws.Range("K106").Value = 220005000
$ws.Range("L106").Value = 4852.6665
$ws.Range("M106").Value = -220004369
$ws.Range("N106").Value = -6114.6665
$ws.Range("H112").Value = 6877.1
$ws.Range("J112").Value = 7983.875
$ws.Range("L112").Value = 23951.625
$ws.Range("N112").Value = -26167.625
$ws.Range("H122").Value = 889.4211
$ws.Range("I122").Value = 844.3889
$ws.Range("K122").Value = 2533.1667
$ws.Range("M122").Value = -83.16670000000022
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H138").Value = 2317.8838
$ws.Range("J138").Value = 2337.08
$ws.Range("L138").Value = 7011.24
$ws.Range("N138").Value = -17291.24

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 307.5
$ws.Range("I4").Value = 320
$ws.Range("J4").Value = 295
$ws.Range("K4").Value = 320
$ws.Range("L4").Value = 295
$ws.Range("M4").Value = -204
$ws.Range("N4").Value = -527
$ws.Range("H32").Value = 9522.357
$ws.Range("I32").Value = 6091.2925
$ws.Range("J32").Value = 18900.6
$ws.Range("K32").Value = 6091.2925
$ws.Range("L32").Value = 18900.6
$ws.Range("M32").Value = -5804.2925
$ws.Range("N32").Value = -19474.6
$ws.Range("H74").Value = 5732.3687
$ws.Range("I74").Value = 2426.7144
$ws.Range("K74").Value = 2426.7144
$ws.Range("M74").Value = -1552.7144
$ws.Range("H77").Value = 5732.3687
$ws.Range("I77").Value = 2426.7144
$ws.Range("K77").Value = 12133.572
$ws.Range("M77").Value = -7765.572
$ws.Range("H97").Value = 1046.5294
$ws.Range("I97").Value = 946.2
$ws.Range("J97").Value = 1799
$ws.Range("K97").Value = 946.2
$ws.Range("L97").Value = 1799
$ws.Range("M97").Value = -450.2
$ws.Range("N97").Value = -2791
$ws.Range("H122").Value = 2028.4872
$ws.Range("I122").Value = 1706.7812
$ws.Range("J122").Value = 3499.1428
$ws.Range("K122").Value = 5120.3436
$ws.Range("L122").Value = 10497.4284
$ws.Range("M122").Value = -2670.3436
$ws.Range("N122").Value = -15397.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2832.75
$ws.Range("I94").Value = 3110.5715
$ws.Range("K94").Value = 3110.5715
$ws.Range("M94").Value = -2659.5715
$ws.Range("H99").Value = 1203.2
$ws.Range("J99").Value = 1328.5
$ws.Range("L99").Value = 1328.5
$ws.Range("N99").Value = -4324.5
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 71851.5
$ws.Range("J9").Value = 71851.5
$ws.Range("L9").Value = 71851.5
$ws.Range("N9").Value = -72187.5
$ws.Range("H16").Value = 1208.8334
$ws.Range("I16").Value = 1050.8
$ws.Range("J16").Value = 1999
$ws.Range("K16").Value = 1050.8
$ws.Range("L16").Value = 1999
$ws.Range("M16").Value = -763.8
$ws.Range("N16").Value = -2573
$ws.Range("H31").Value = 12798.454
$ws.Range("I31").Value = 7009.4614
$ws.Range("K31").Value = 7009.4614
$ws.Range("M31").Value = -6714.4614
$ws.Range("H34").Value = 12798.454
$ws.Range("I34").Value = 7009.4614
$ws.Range("K34").Value = 7009.4614
$ws.Range("M34").Value = -6807.4614
$ws.Range("H113").Value = 1208.8334
$ws.Range("I113").Value = 1050.8
$ws.Range("J113").Value = 1999
$ws.Range("K113").Value = 1050.8
$ws.Range("L113").Value = 1999
$ws.Range("M113").Value = 1119.2
$ws.Range("N113").Value = -6339
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 6373.8423
$ws.Range("I132").Value = 2864.5715
$ws.Range("K132").Value = 8593.7145
$ws.Range("M132").Value = -6063.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 750
$ws.Range("I55").Value = 750
$ws.Range("K55").Value = 2250
$ws.Range("M55").Value = -2073
$ws.Range("H131").Value = 4606.095
$ws.Range("I131").Value = 2784.5
$ws.Range("J131").Value = 5334.7334
$ws.Range("K131").Value = 8353.5
$ws.Range("L131").Value = 16004.2002
$ws.Range("M131").Value = -3313.5
$ws.Range("N131").Value = -26084.2002
$ws.Range("H134").Value = 2079.5386
$ws.Range("I134").Value = 2079.5386
$ws.Range("K134").Value = 6238.6158
$ws.Range("M134").Value = -1168.6158
$ws.Range("H141").Value = 3957.3845
$ws.Range("I141").Value = 3949.7273
$ws.Range("J141").Value = 3999.5
$ws.Range("K141").Value = 11849.1819
$ws.Range("L141").Value = 11998.5
$ws.Range("M141").Value = -6669.1819
$ws.Range("N141").Value = -22358.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 759.25
$ws.Range("I97").Value = 691.9091
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 691.9091
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -195.9091
$ws.Range("N97").Value = -2492
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H113").Value = 4501.8335
$ws.Range("I113").Value = 2670.3333
$ws.Range("J113").Value = 6333.3335
$ws.Range("K113").Value = 2670.3333
$ws.Range("L113").Value = 6333.3335
$ws.Range("M113").Value = -500.3332999999998
$ws.Range("N113").Value = -10673.3335
$ws.Range("H138").Value = 86000
$ws.Range("J138").Value = 86000
$ws.Range("L138").Value = 86000
$ws.Range("N138").Value = -96280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 60000
$ws.Range("J38").Value = 60000
$ws.Range("L38").Value = 60000
$ws.Range("N38").Value = -60820
$ws.Range("H46").Value = 2883.5
$ws.Range("J46").Value = 3586.6924
$ws.Range("L46").Value = 3586.6924
$ws.Range("N46").Value = -3962.6924
$ws.Range("H61").Value = 7661
$ws.Range("I61").Value = 1573.75
$ws.Range("J61").Value = 13748.25
$ws.Range("K61").Value = 1573.75
$ws.Range("L61").Value = 13748.25
$ws.Range("M61").Value = -1371.75
$ws.Range("N61").Value = -14152.25
$ws.Range("H113").Value = 7661
$ws.Range("I113").Value = 1573.75
$ws.Range("J113").Value = 13748.25
$ws.Range("K113").Value = 1573.75
$ws.Range("L113").Value = 13748.25
$ws.Range("M113").Value = 596.25
$ws.Range("N113").Value = -18088.25
$ws.Range("H130").Value = 29999.143
$ws.Range("J130").Value = 29999.143
$ws.Range("L130").Value = 29999.143
$ws.Range("N130").Value = -40039.143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 731.6
$ws.Range("I100").Value = 356.33334
$ws.Range("J100").Value = 1294.5
$ws.Range("K100").Value = 712.66668
$ws.Range("L100").Value = 2589
$ws.Range("M100").Value = -171.66668
$ws.Range("N100").Value = -3671
$ws.Range("H107").Value = 1132.1224
$ws.Range("I107").Value = 1196.0322
$ws.Range("K107").Value = 3588.0966
$ws.Range("M107").Value = -1668.0966
$ws.Range("H136").Value = 6269.314
$ws.Range("I136").Value = 3307.6667
$ws.Range("J136").Value = 13377.267
$ws.Range("K136").Value = 9923.000100000001
$ws.Range("L136").Value = 40131.801
$ws.Range("M136").Value = -7373.000100000001
$ws.Range("N136").Value = -45231.801

Write-Output "Applied changes: sets=243 clears=5"
